$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Days")

# --- Add 10 new rows (86-95) continuing the existing Day/Date sequence ---
# Row 85 currently holds day 84 / serial date 43549 (2019-03-25); the new
# rows continue with days 85-94 and serial dates 43550-43559.

# Copy the formatting (number format, font, borders, alignment) from the
# last existing data row so the new rows reuse the same cell styles rather
# than creating new ones.
$ws.Range("A85").Copy()
$ws.Range("A86:A95").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B85").Copy()
$ws.Range("B86:B95").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

$startRow = 86
$startDay = 85
$startSerial = 43550

for ($i = 0; $i -lt 10; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $startDay + $i
    $ws.Cells.Item($r, 2).Value = $startSerial + $i
    $ws.Rows.Item($r).RowHeight = 13.8
}

# --- Update the active selection to C83 (was C79) ---
$ws.Range("C83").Select()

# --- Adjust the workbook tab-bar split ratio ---
$excel.ActiveWindow.TabRatio = 500
